# Update the Green Mountain salads worksheet:
#  - Replace the placeholder salad rows with the completed nutrition data
#    (Cobb Salad / Caesar Salad / House Salad) including allergens, local
#    ingredients, diet codes and nutrition-label keys.
#  - Re-centre the two "local ingredients" highlight cells.
#  - Update the current selection to the data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Cobb Salad -----------------------------------------------------
$ws.Range("A2").Value = "Cobb Salad"
$ws.Range("B2").Value = " Egg / Guacamole / Bacon / Cheddar / Tomato / Cucumber / Romaine / Ranch Dressing"
$ws.Range("C2").Value = "Milk, eggs."
$ws.Range("D2").Value = "Fresh Start Cucumber, Fresh Start Tomato, Fresh Start Romaine"
$ws.Range("E2").Value = "BC"
$ws.Range("F2").Value = "Cobb_Salad"

# --- Row 3: Caesar Salad ----------------------------------------------------
$ws.Range("A3").Value = "Caesar Salad"
$ws.Range("B3").Value = "Crispy Chicken / Bacon / Romaine / Parmesan / Caesar Dressing"
$ws.Range("C3").Value = "Milk, eggs, fish."
$ws.Range("D3").Value = "Fresh Start Romaine, Castle Cheese Parmesan"
$ws.Range("E3").Value = "BC"
$ws.Range("F3").Value = "Caesar_Salad"

# --- Row 4: House Salad ------------------------------------------------------
$ws.Range("A4").Value = "House Salad"
$ws.Range("B4").Value = "Carrots / Tomato / Cucumber / Romaine / Balsamic Dressing"
$ws.Range("C4").Value = "Dressing contains sulphites."
$ws.Range("D4").Value = "Fresh Start Cucumber, Fresh Start Tomato, Fresh Start Romaine"
$ws.Range("E4").Value = "BC, VEG, VGN, GF DF"
$ws.Range("F4").Value = "House_Salad"

# Clear the old "wrap text" formatting that used to sit on B3/C3 and instead
# centre (horizontally + vertically) the two local-ingredients callout cells.
$ws.Range("B3").WrapText = $false
$ws.Range("C3").WrapText = $false

$b2 = $ws.Range("B2")
$b2.HorizontalAlignment = -4108
$b2.VerticalAlignment = -4108

$b4 = $ws.Range("B4")
$b4.HorizontalAlignment = -4108
$b4.VerticalAlignment = -4108

# Select the newly completed data rows (whole rows 2-4), as left by the author.
$ws.Activate()
$ws.Rows("2:4").Select()
